# 添加 2022-Q3 数据 (feat: add 2022-Q3 data)
#
# 1. Create a new "2022-Q3" sheet (as a copy of "2022-Q2" so it inherits the
#    same cell styling/borders), positioned right before "2022-Q2", and fill
#    it in with the three 2022-Q3 fund rows.
# 2. Insert a new row at the top of the "总计" (totals) summary sheet with
#    the 2022-Q3 aggregate figures, pushing the existing quarters down.
# 3. Restore the originally-selected/active sheet ("2021-Q1", the last tab).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: build the "2022-Q3" worksheet
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)                      # duplicate placed immediately before "2022-Q2"
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# "2022-Q2" has 6 fund rows (rows 2-7); 2022-Q3 only has 3, drop the extra 3.
$q3.Rows.Item(5).Delete()
$q3.Rows.Item(5).Delete()
$q3.Rows.Item(5).Delete()

# Columns B-G hold text (fund code / name / percentages kept as strings);
# force text formatting before assigning so leading zeros & decimal text
# are preserved instead of being coerced to numbers.
$q3.Range("B2:G4").NumberFormat = "@"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "006039"
$q3.Range("C2").Value = "国富估值优势混合"
$q3.Range("D2").Value = "6.13"
$q3.Range("E2").Value = "81.55"
$q3.Range("F2").Value = "2.95"
$q3.Range("G2").Value = "0.1808"
$q3.Range("H2").Value = 5

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "014339"
$q3.Range("C3").Value = "长江智能制造混合A"
$q3.Range("D3").Value = "2.42"
$q3.Range("E3").Value = "75.36"
$q3.Range("F3").Value = "3.23"
$q3.Range("G3").Value = "0.0782"
$q3.Range("H3").Value = 6

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "014340"
$q3.Range("C4").Value = "长江智能制造混合C"
$q3.Range("D4").Value = "0.11"
$q3.Range("E4").Value = "75.36"
$q3.Range("F4").Value = "3.23"
$q3.Range("G4").Value = "0.0036"
$q3.Range("H4").Value = 6

# ---------------------------------------------------------------------------
# Step 2: update the "总计" summary sheet
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.26

# Re-number the index column (0,1,2,...) now that a row was inserted.
for ($i = 0; $i -le 5; $i++) {
    $total.Cells.Item($i + 2, 1).Value = $i
}

# ---------------------------------------------------------------------------
# Step 3: restore the originally active/selected sheet
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
$wb.Worksheets.Item("2021-Q1").Range("A1").Select()
